$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, shifting existing rows 19-44 down to 20-45.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new weekly price record.
$ws.Range("A19").Value = 10
$ws.Range("B19").Value = "Vega Modelo de Temuco"
$ws.Range("C19").Value = "La Araucanía"
$ws.Range("D19").Value = 44546
$ws.Range("E19").Value = 9
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100103
$ws.Range("H19").Value = "Frutos de hueso (carozo)"
$ws.Range("I19").Value = 100103003
$ws.Range("J19").Value = "Damasco"
$ws.Range("K19").Value = "Castle Brite"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 65
$ws.Range("N19").Value = 18000
$ws.Range("O19").Value = 18000
$ws.Range("P19").Value = 18000
$ws.Range("Q19").Value = "`$/bandeja 18 kilos"
$ws.Range("R19").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S19").Value = 1000
$ws.Range("T19").Value = 18
